# Add a title block ("LvlUpFitness: " / "Feasibility Study : ") as two new
# paragraphs at the very beginning of the document, before the existing
# first paragraph.

$d = $word.ActiveDocument

$wordMl = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$xml = '<w:p xmlns:w="' + $wordMl + '">' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>LvlUpFitness</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
       '</w:p>' +
       '<w:p xmlns:w="' + $wordMl + '">' +
         '<w:r><w:t>Feasibility</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> Study : </w:t></w:r>' +
       '</w:p>'

# Insert the two new paragraphs at the very start of the document (position 0).
$start = $d.Range(0, 0)
$start.InsertXML($xml)

Write-Output "Inserted title paragraphs"
